# Natmi following Dr Hou advice
# Rewrites the LR-pairs table (Wnt4-Fzd6) with the recomputed values
# covering all Sending-cluster x Target-cluster combinations (ECs, FAPs, M1, M2, Neutro, sCs) x (ECs, FAPs, sCs).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs","Wnt4","Fzd6","ECs",3,1,0.8155043333333333,2.446513,0.06451520947540776,0.06451520947540779,3,1,13.877148,41.631444,0.9551716529386821,0.9551716529386822,11.316874328308,101.851868954772,0.06162309927431055,0.06162309927431059),
    @("ECs","Wnt4","Fzd6","FAPs",3,1,0.8155043333333333,2.446513,0.06451520947540776,0.06451520947540779,2,0.6666666666666666,0.4259926666666667,1.277978,0.02932130719941569,0.02932130719941569,0.3473988656348889,3.126589790714,0.001891670276063085,0.001891670276063086),
    @("ECs","Wnt4","Fzd6","sCs",3,1,0.8155043333333333,2.446513,0.06451520947540776,0.06451520947540779,3,1,0.225293,0.675879,0.01550703986190206,0.01550703986190207,0.1837274177696666,1.653546759927,0.00100043992503411,0.00100043992503411),
    @("FAPs","Wnt4","Fzd6","ECs",3,1,1.070199666666667,3.210599,0.08466436394433004,0.08466436394433006,3,1,13.877148,41.631444,0.9551716529386821,0.9551716529386822,14.851319163884,133.661872474956,0.08086900045370789,0.08086900045370791),
    @("FAPs","Wnt4","Fzd6","FAPs",3,1,1.070199666666667,3.210599,0.08466436394433004,0.08466436394433006,2,0.6666666666666666,0.4259926666666667,1.277978,0.02932130719941569,0.02932130719941569,0.4558972098691111,4.103074888822,0.002482469824054835,0.002482469824054835),
    @("FAPs","Wnt4","Fzd6","sCs",3,1,1.070199666666667,3.210599,0.08466436394433004,0.08466436394433006,3,1,0.225293,0.675879,0.01550703986190206,0.01550703986190207,0.2411084935023333,2.169976441521,0.00131289366656731,0.00131289366656731),
    @("M1","Wnt4","Fzd6","ECs",3,1,1.792556333333333,5.377669,0.1418105859336969,0.1418105859336969,3,1,13.877148,41.631444,0.9551716529386821,0.9551716529386822,24.875569536004,223.880125824036,0.1354534517704923,0.1354534517704923),
    @("M1","Wnt4","Fzd6","FAPs",3,1,1.792556333333333,5.377669,0.1418105859336969,0.1418105859336969,2,0.6666666666666666,0.4259926666666667,1.277978,0.02932130719941569,0.02932130719941569,0.7636158525868889,6.872542673282,0.004158071754291065,0.004158071754291066),
    @("M1","Wnt4","Fzd6","sCs",3,1,1.792556333333333,5.377669,0.1418105859336969,0.1418105859336969,3,1,0.225293,0.675879,0.01550703986190206,0.01550703986190207,0.4038503940056667,3.634653546051,0.002199062408913526,0.002199062408913527),
    @("M2","Wnt4","Fzd6","ECs",3,1,1.782737666666667,5.348212999999999,0.1410338232472499,0.1410338232472499,3,1,13.877148,41.631444,0.9551716529386821,0.9551716529386822,24.739314445508,222.653830009572,0.1347115100713376,0.1347115100713376),
    @("M2","Wnt4","Fzd6","FAPs",3,1,1.782737666666667,5.348212999999999,0.1410338232472499,0.1410338232472499,2,0.6666666666666666,0.4259926666666667,1.277978,0.02932130719941569,0.02932130719941569,0.7594331725904444,6.834898553314,0.004135296056940707,0.004135296056940708),
    @("M2","Wnt4","Fzd6","sCs",3,1,1.782737666666667,5.348212999999999,0.1410338232472499,0.1410338232472499,3,1,0.225293,0.675879,0.01550703986190206,0.01550703986190207,0.4016383171363333,3.614744854227,0.002187017118971553,0.002187017118971554),
    @("Neutro","Wnt4","Fzd6","ECs",3,1,6.337404666666667,19.012214,0.5013572250796461,0.5013572250796462,3,1,13.877148,41.631444,0.9551716529386821,0.9551716529386822,87.94510249522401,791.505922457016,0.4788822093920764,0.4788822093920766),
    @("Neutro","Wnt4","Fzd6","FAPs",3,1,6.337404666666667,19.012214,0.5013572250796461,0.5013572250796462,2,0.6666666666666666,0.4259926666666667,1.277978,0.02932130719941569,0.02932130719941569,2.699687913699111,24.297191223292,0.0147004492132069,0.0147004492132069),
    @("Neutro","Wnt4","Fzd6","sCs",3,1,6.337404666666667,19.012214,0.5013572250796461,0.5013572250796462,3,1,0.225293,0.675879,0.01550703986190206,0.01550703986190207,1.427772909567333,12.849956186106,0.007774566474362677,0.007774566474362679),
    @("sCs","Wnt4","Fzd6","ECs",3,1,0.8420946666666667,2.526284,0.06661879231966927,0.06661879231966929,3,1,13.877148,41.631444,0.9551716529386821,0.9551716529386822,11.685872319344,105.172850874096,0.06363238197675727,0.0636323819767573),
    @("sCs","Wnt4","Fzd6","FAPs",3,1,0.8420946666666667,2.526284,0.06661879231966927,0.06661879231966929,2,0.6666666666666666,0.4259926666666667,1.277978,0.02932130719941569,0.02932130719941569,0.3587261526391111,3.228535373752,0.001953350074859097,0.001953350074859098),
    @("sCs","Wnt4","Fzd6","sCs",3,1,0.8420946666666667,2.526284,0.06661879231966927,0.06661879231966929,3,1,0.225293,0.675879,0.01550703986190206,0.01550703986190207,0.1897180337373333,1.707462303636,0.001033060268052886,0.001033060268052887)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($startRow + $i, $j + 1).Value = $row[$j]
    }
}
